$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "MLP-deep"
$ws.Cells.Item(2, 2).Value = 0.6153999999999999
$ws.Cells.Item(2, 3).Value = 0.6153999999999999
$ws.Cells.Item(2, 4).Value = 0.9451000000000001
$ws.Cells.Item(2, 5).Value = 0.5367
$ws.Cells.Item(2, 6).Value = 0.5779
$ws.Cells.Item(2, 7).Value = 0.7018
$ws.Cells.Item(2, 8).Value = 0.7010999999999999
$ws.Cells.Item(3, 1).Value = "MLP 128"
$ws.Cells.Item(3, 2).Value = 0.5412
$ws.Cells.Item(3, 3).Value = 0.5412
$ws.Cells.Item(3, 4).Value = 0.9409
$ws.Cells.Item(3, 5).Value = 0.594
$ws.Cells.Item(3, 6).Value = 0.6035
$ws.Cells.Item(3, 7).Value = 0.6894
$ws.Cells.Item(3, 8).Value = 0.6879
$ws.Cells.Item(4, 1).Value = "MLP 16"
$ws.Cells.Item(4, 2).Value = 0.5302
$ws.Cells.Item(4, 3).Value = 0.5302
$ws.Cells.Item(4, 4).Value = 0.9382
$ws.Cells.Item(4, 5).Value = 0.6158
$ws.Cells.Item(4, 6).Value = 0.6322
$ws.Cells.Item(4, 7).Value = 0.6731
$ws.Cells.Item(4, 8).Value = 0.6731
$ws.Cells.Item(5, 1).Value = "MLP 32"
$ws.Cells.Item(5, 2).Value = 0.5371
$ws.Cells.Item(5, 3).Value = 0.5371
$ws.Cells.Item(5, 4).Value = 0.9354
$ws.Cells.Item(5, 5).Value = 0.6151
$ws.Cells.Item(5, 6).Value = 0.6466
$ws.Cells.Item(5, 7).Value = 0.6666
$ws.Cells.Item(5, 8).Value = 0.6656
$ws.Cells.Item(6, 1).Value = "MLP 64"
$ws.Cells.Item(6, 2).Value = 0.5646
$ws.Cells.Item(6, 3).Value = 0.5646
$ws.Cells.Item(6, 4).Value = 0.9354
$ws.Cells.Item(6, 5).Value = 0.5927
$ws.Cells.Item(6, 6).Value = 0.6084000000000001
$ws.Cells.Item(6, 7).Value = 0.6876
$ws.Cells.Item(6, 8).Value = 0.6853
$ws.Cells.Item(7, 1).Value = "SVR rbf"
$ws.Cells.Item(7, 2).Value = 0.5206
$ws.Cells.Item(7, 3).Value = 0.5206
$ws.Cells.Item(7, 4).Value = 0.9258
$ws.Cells.Item(7, 5).Value = 0.6096
$ws.Cells.Item(7, 6).Value = 0.639
$ws.Cells.Item(7, 7).Value = 0.6697
$ws.Cells.Item(7, 8).Value = 0.6695
$ws.Cells.Item(8, 1).Value = "LinearRegression"
$ws.Cells.Item(8, 2).Value = 0.4533
$ws.Cells.Item(8, 3).Value = 0.4533
$ws.Cells.Item(8, 4).Value = 0.9147999999999999
$ws.Cells.Item(8, 5).Value = 0.7194
$ws.Cells.Item(8, 6).Value = 0.8163
$ws.Cells.Item(8, 7).Value = 0.5784
$ws.Cells.Item(8, 8).Value = 0.5778
$ws.Cells.Item(9, 1).Value = "LinearSVR"
$ws.Cells.Item(9, 2).Value = 0.4808
$ws.Cells.Item(9, 3).Value = 0.4808
$ws.Cells.Item(9, 4).Value = 0.9121
$ws.Cells.Item(9, 5).Value = 0.7091
$ws.Cells.Item(9, 6).Value = 0.8129999999999999
$ws.Cells.Item(9, 7).Value = 0.5796
$ws.Cells.Item(9, 8).Value = 0.5796
$ws.Cells.Item(10, 1).Value = "RidgeCV"
$ws.Cells.Item(10, 2).Value = 0.456
$ws.Cells.Item(10, 3).Value = 0.456
$ws.Cells.Item(10, 4).Value = 0.9121
$ws.Cells.Item(10, 5).Value = 0.7186
$ws.Cells.Item(10, 6).Value = 0.8117
$ws.Cells.Item(10, 7).Value = 0.5808
$ws.Cells.Item(10, 8).Value = 0.5802
$ws.Cells.Item(11, 1).Value = "GradientBoostingRegressor"
$ws.Cells.Item(11, 2).Value = 0.4973
$ws.Cells.Item(11, 3).Value = 0.4973
$ws.Cells.Item(11, 4).Value = 0.9107
$ws.Cells.Item(11, 5).Value = 0.6563
$ws.Cells.Item(11, 6).Value = 0.7339
$ws.Cells.Item(11, 7).Value = 0.6205000000000001
$ws.Cells.Item(11, 8).Value = 0.6205000000000001
$ws.Cells.Item(12, 1).Value = "LassoCV"
$ws.Cells.Item(12, 2).Value = 0.4615
$ws.Cells.Item(12, 3).Value = 0.4615
$ws.Cells.Item(12, 4).Value = 0.9066
$ws.Cells.Item(12, 5).Value = 0.717
$ws.Cells.Item(12, 6).Value = 0.8115
$ws.Cells.Item(12, 7).Value = 0.5809
$ws.Cells.Item(12, 8).Value = 0.5803
$ws.Cells.Item(13, 1).Value = "SVR sigmoid"
$ws.Cells.Item(13, 2).Value = 0.3187
$ws.Cells.Item(13, 3).Value = 0.3187
$ws.Cells.Item(13, 4).Value = 0.8269
$ws.Cells.Item(13, 5).Value = 0.9006
$ws.Cells.Item(13, 6).Value = 1.2235
$ws.Cells.Item(13, 7).Value = 0.3845
$ws.Cells.Item(13, 8).Value = 0.3673
$ws.Cells.Item(14, 1).Value = "SVR poly"
$ws.Cells.Item(14, 2).Value = 0.283
$ws.Cells.Item(14, 3).Value = 0.283
$ws.Cells.Item(14, 4).Value = 0.7706
$ws.Cells.Item(14, 5).Value = 0.9792999999999999
$ws.Cells.Item(14, 6).Value = 2.1647
$ws.Cells.Item(14, 7).Value = -0.1195
$ws.Cells.Item(14, 8).Value = -0.1195
